$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 29: Total_tested (column R) changes from 0 to 1191
$ws.Range("R29").Value = 1191

# Append new row 30 with the latest TCHD data
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 70
$ws.Range("Q30").Value = 996
$ws.Range("R30").Value = 1197
$ws.Range("S30").Value = 0
$ws.Range("T30").Value = 0
$ws.Range("U30").Value = 0
$ws.Range("V30").Value = 2

# Match the style used for column A on preceding data rows (bordered/bold/centered)
$ws.Range("A29").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
